$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @{
  "G2" = 3.8
  "I2" = 2.22
  "J2" = 3.6
  "K2" = 3.65
  "V2" = 1.82
  "W2" = 1.35
  "AK2" = 48

  "F3" = 1.12
  "G3" = 2
  "H3" = 2.4
  "I3" = 36
  "J3" = 2.1
  "K3" = 12
  "L3" = 1.2
  "R3" = 1.33
  "S3" = 1.8
  "V3" = 1.02
  "W3" = 2

  "G4" = 5.8
  "I4" = 1.89
  "J4" = 3.45
  "Q4" = 1.65
  "T4" = 1.67
  "U4" = 2.22
  "V4" = 2.12
  "W4" = 1.23
  "AO4" = 1000

  "I5" = 2.18
  "K5" = 3.95
  "V5" = 1.84

  "H6" = 1.97
  "K6" = 4.2

  "G7" = 9.4
  "I7" = 1.49
  "K7" = 5.6

  "G8" = 3.55
  "H8" = 2.04
  "K8" = 4.6
  "N8" = 5.9
  "U8" = 2.66
  "W8" = 1.4
  "Y8" = 19
  "Z8" = 22
  "AA8" = 34
  "AD8" = 12
  "AE8" = 20
  "AF8" = 36
  "AG8" = 18.5
  "AI8" = 32
  "AL8" = 34
  "AN8" = 20
  "AO8" = 11.5

  "K9" = 4.3
  "AN9" = 13.5

  "G10" = 6.4
  "H10" = 1.63
  "J10" = 4.1
  "K10" = 4.4
  "W10" = 1.19
  "AE10" = 1000
  "AG10" = 25
  "AL10" = 85

  "AF11" = 170
  "AJ11" = 500
  "AK11" = 200
  "AL11" = 140
}

foreach ($addr in $updates.Keys) {
  $ws.Range($addr).Value = $updates[$addr]
}
